# Generate Report for Handoff
#
# The localization-status report rotates the three rows that track the
# files "ed11fd94...md", "2a303c20...md" and "8e2ffb84...md":
#   - 2a303c20 and 8e2ffb84 (already "Ready for handoff") move up one row
#   - ed11fd94 is freshly handed off, moves to the bottom of that block,
#     its status becomes "Ready for handoff", and it gets a new
#     "Latest Handoff Datetime" on the zh-cn / de-de sheets.
# This is applied identically (by row position) on the Overview, zh-cn
# and de-de sheets. Hyperlink targets (the rIds / URLs) stay attached to
# the same row position, only the displayed text changes - so we update
# both the cell values and each hyperlink's TextToDisplay.

$wb = $excel.ActiveWorkbook

function Set-CellAndLink {
    param(
        $ws,
        [int]$row,
        [int]$col,
        [string]$text
    )

    $ws.Cells.Item($row, $col).Value = $text

    foreach ($hl in $ws.Hyperlinks) {
        $hlRow = $hl.Range.Row()
        $hlCol = $hl.Range.Column()
        if ($hlRow -eq $row -and $hlCol -eq $col) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet: columns A (File Name), B (zh-cn status), C (de-de status)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellAndLink $wsOverview 7 1 "2a303c20-f4a1-45c8-962c-f50d0750d8fd.md"
$wsOverview.Cells.Item(7, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item(7, 3).Value = "Ready for handoff"

Set-CellAndLink $wsOverview 8 1 "8e2ffb84-5204-4ff0-a15a-ed1a502b4430.md"
$wsOverview.Cells.Item(8, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item(8, 3).Value = "Ready for handoff"

Set-CellAndLink $wsOverview 9 1 "ed11fd94-8665-446d-a6d1-4cdc3d012e2c.md"
$wsOverview.Cells.Item(9, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item(9, 3).Value = "Ready for handoff"

# ---------------------------------------------------------------------
# zh-cn sheet: A=Source File Name, B=Status, C=Latest Handoff File,
#              D=Latest Handoff Datetime
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellAndLink $wsZhCn 7 1 "2a303c20-f4a1-45c8-962c-f50d0750d8fd.md"
$wsZhCn.Cells.Item(7, 2).Value = "Ready for handoff"
Set-CellAndLink $wsZhCn 7 3 "2a303c20-f4a1-45c8-962c-f50d0750d8fd.b923017c05e9db8a2cc6376f21bd5bf241c674b0.zh-cn.xlf"
$wsZhCn.Cells.Item(7, 4).Value = "2016-03-10 08:59:27"

Set-CellAndLink $wsZhCn 8 1 "8e2ffb84-5204-4ff0-a15a-ed1a502b4430.md"
$wsZhCn.Cells.Item(8, 2).Value = "Ready for handoff"
Set-CellAndLink $wsZhCn 8 3 "8e2ffb84-5204-4ff0-a15a-ed1a502b4430.f7b80f82519d24f60afaf83ba6163fc878b9e900.zh-cn.xlf"
$wsZhCn.Cells.Item(8, 4).Value = "2016-03-10 09:07:06"

Set-CellAndLink $wsZhCn 9 1 "ed11fd94-8665-446d-a6d1-4cdc3d012e2c.md"
$wsZhCn.Cells.Item(9, 2).Value = "Ready for handoff"
Set-CellAndLink $wsZhCn 9 3 "ed11fd94-8665-446d-a6d1-4cdc3d012e2c.1e9e0e762e630aaf9711ba42142d001945c6c503.zh-cn.xlf"
$wsZhCn.Cells.Item(9, 4).Value = "2016-03-10 09:15:28"

# ---------------------------------------------------------------------
# de-de sheet: A=Source File Name, B=Status, C=Latest Handoff File,
#              D=Latest Handoff Datetime
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellAndLink $wsDeDe 7 1 "2a303c20-f4a1-45c8-962c-f50d0750d8fd.md"
$wsDeDe.Cells.Item(7, 2).Value = "Ready for handoff"
Set-CellAndLink $wsDeDe 7 3 "2a303c20-f4a1-45c8-962c-f50d0750d8fd.b923017c05e9db8a2cc6376f21bd5bf241c674b0.de-de.xlf"
$wsDeDe.Cells.Item(7, 4).Value = "2016-03-10 08:59:33"

Set-CellAndLink $wsDeDe 8 1 "8e2ffb84-5204-4ff0-a15a-ed1a502b4430.md"
$wsDeDe.Cells.Item(8, 2).Value = "Ready for handoff"
Set-CellAndLink $wsDeDe 8 3 "8e2ffb84-5204-4ff0-a15a-ed1a502b4430.f7b80f82519d24f60afaf83ba6163fc878b9e900.de-de.xlf"
$wsDeDe.Cells.Item(8, 4).Value = "2016-03-10 09:07:14"

Set-CellAndLink $wsDeDe 9 1 "ed11fd94-8665-446d-a6d1-4cdc3d012e2c.md"
$wsDeDe.Cells.Item(9, 2).Value = "Ready for handoff"
Set-CellAndLink $wsDeDe 9 3 "ed11fd94-8665-446d-a6d1-4cdc3d012e2c.1e9e0e762e630aaf9711ba42142d001945c6c503.de-de.xlf"
$wsDeDe.Cells.Item(9, 4).Value = "2016-03-10 09:15:35"

Write-Output "Done updating handoff rotation."
